$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 updates
$ws.Range("I11").Value = 18999
$ws.Range("K11").Value = 103695
$ws.Range("R11").Value = 19000.8999
$ws.Range("T11").Value = 103705.3695

# Row 12 updates
$ws.Range("I12").Value = 19289
$ws.Range("K12").Value = 107380
$ws.Range("R12").Value = 19290.9289
$ws.Range("T12").Value = 107390.738
